$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(17.96285899517179,8.714392610023914,8.018325215040209,12.81274015110869,36.38381813520179,10.0589417986614,10.95819561798781,16.91670573382469,19.49425914206328,28.01996576738953),
    @(17.62360461244337,8.516187486938058,8.023061227770761,12.84045685644303,36.45166563294448,10.07565284077265,10.96114853881381,16.84883679745139,19.55327155283013,28.08663974596839),
    @(17.41491989629298,8.391035433927787,8.026676513517501,12.85848449738514,36.50140641228868,10.08643842088905,10.96409290700717,16.80930443695761,19.59135150766563,28.13351603022231),
    @(17.32989696240397,8.339211677544784,8.028328138739752,12.86608535368003,36.52370500449121,10.0909660190077,10.96557811437714,16.79374380987196,19.60733484165292,28.15410887254677),
    @(17.31578340534548,8.33055802827525,8.028613176962972,12.86736285587661,36.52753009212098,10.09172583047263,10.96584199153617,16.79119346037586,19.61001700996794,28.15761822581502),
    @(17.41377302357473,8.390339792512771,8.026698065060293,12.85858597413633,36.50169892931503,10.08649894511634,10.96411178051005,16.80909234319097,19.59156517830712,28.13378772196553),
    @(17.84603063436421,8.646794423067822,8.01981168305803,12.82208772043336,36.4055325485127,10.06459505047874,10.95897953747855,16.8928673387155,19.51422420345074,28.04172085213545),
    @(18.68580417460907,9.120460861147556,8.011898067654082,12.75849692486836,36.28121590091604,10.02578884954097,10.95785121921217,17.07361231522736,19.37715044432072,27.90841277700157),
    @(19.291719711254,9.448374094535081,8.009462304038776,12.71660455820026,36.2292168087527,9.999781073023376,10.96241448552366,17.21572265118291,19.28525864185478,27.83941493266225),
    @(19.56372346085555,9.592765311697507,8.009081096547579,12.6985868344056,36.21412465905376,9.988487638915849,10.96564845905668,17.28223252467275,19.2453523562526,27.81433688698502),
    @(19.6661055612616,9.646723727497289,8.009040651571016,12.69191283547256,36.2096416162268,9.984288019403913,10.967038412363,17.3076715634998,19.23051226962888,27.80574936777033),
    @(19.6440848088027,9.635135298810464,8.009044750607591,12.69334358560981,36.2105523200812,9.985189065415994,10.96673172853455,17.30218178876152,19.23369628672439,27.80755839236491),
    @(19.57215946040269,9.597219096484064,8.009075690161833,12.69803477890029,36.21373114282476,9.988140593327259,10.96575950532924,17.28432042096259,19.24412601692244,27.81361216045659),
    @(19.5280195812181,9.573899727104026,8.009108155274767,12.70092764861949,36.2158387207791,9.989958501312914,10.96518548261167,17.273412347463,19.25054986180282,27.81743869083291),
    @(19.27386209191125,9.438838816352794,8.009501791143615,12.71780292710944,36.23037549471317,10.00052991389302,10.96222633927998,17.21141231557402,19.28790466876897,27.84118095688666),
    @(19.11694359877237,9.354735670428187,8.009929008040507,12.72842117959557,36.24148713787596,10.00715258247372,10.96070677496988,17.1738441449856,19.31130542876007,27.85736337945674),
    @(19.02635204333027,9.305913779604804,8.010243182677131,12.7346263750777,36.24868420036368,10.01101239178612,10.95994184624022,17.15241220314711,19.32494342608207,27.86726487181525),
    @(18.99562462771779,9.289307613175655,8.010361331449991,12.73674417105679,36.25125938932289,10.01232796266106,10.95970162603849,17.14518643425214,19.32959170733946,27.87071927384602),
    @(19.13368336026589,9.363735201353546,8.009876450014389,12.7272807244806,36.24022086525996,10.00644235110749,10.96085725399264,17.17782520247736,19.30879591157635,27.85557926439596),
    @(19.59330325229802,9.608375772799594,8.009063787329286,12.69665282421984,36.21276400597555,9.987271572183557,10.96604059418205,17.28955998519716,19.24105519062017,27.81180934027667),
    @(19.89003723522556,9.764056180043996,8.009138034976843,12.67750351004437,36.20200043460335,9.97519077390297,10.97039119768078,17.36405385631076,19.19836508535598,27.78850202362442),
    @(19.73202995252267,9.681361481443796,8.009043222692309,12.6876446306871,36.20708801009873,9.981597608989302,10.96798150286697,17.32416564143833,19.22100515289263,27.80045625373992),
    @(19.12611648091769,9.359667974525779,8.009899997901924,12.72779601077957,36.24079082769672,10.00676328361767,10.96078888378221,17.17602484662837,19.30992988921721,27.85638400018327),
    @(18.4601294075697,8.995696592002199,8.013443396787927,12.77484936061652,36.30794810966361,10.03584557192942,10.95720529110219,17.023025827842,19.41267864351311,27.93940292395218)
)

$cols = @("B","C","D","E","F","J","L","M","N","O")
$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $data[$i][$j]
    }
}
